# Parts list update:
#  - Item 9 (row 11, "TRRS to TRS Adapter" from Amazon) didn't work out for the
#    author's purposes, so it gets moved down into the "EXTRA COSTS" section
#    (row 18) with a note explaining why, and item 9's slot is replaced with
#    the part the author actually ended up using.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move row 11 (A:G) down into row 18, the first open row under EXTRA COSTS.
# Cut leaves the formatting (style indexes) behind on row 11, which is what we
# want since row 11 is about to be populated with a brand new part.
$ws.Range("A11:G11").Cut($ws.Range("A18:G18"))

# Row 18 is now the 5th EXTRA COSTS line item (numbered 0-3 already occupy
# rows 14-17), and the reason it ended up here is different from its original
# note, so fix both up.
$ws.Range("A18").Value = 4
$ws.Range("G18").Value = "Not what I needed, going to make my own"

# Row 11 becomes the new part that replaced the old TRRS to TRS adapter.
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "TRRS female connector mount"
$ws.Range("C11").ClearContents()
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 6.34
$ws.Range("F11").Value = "ShowMeCables"
$ws.Range("G11").Value = "Going to make own TRRS to TRS connector"

# Match the vertical-centered alignment used by the other item rows (the old
# row 11 had only horizontal centering).
$ws.Range("A11").VerticalAlignment = -4108

# Extend the cost total to cover the newly-filled row 18.
$ws.Range("E20").Formula = "=SUM(E2:E18)"

# Column B needs to be a bit wider to fit the new part name.
$ws.Columns("B").ColumnWidth = 25.2

# Leave the selection where the author finished editing.
$ws.Range("E21").Select() | Out-Null
